$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.762.14"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "2.237.92"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "315.72"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").Value = "98.79"
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("D7").Value = "0.568"
$ws.Range("E7").Value = "  -2.97%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  -5.86%  "
$ws.Range("D10").Value = "36.17"
$ws.Range("E10").Value = "  -3.91%  "
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("D12").Value = "7.35"
$ws.Range("E12").Value = "  -5.55%  "
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("D14").Value = "2.569.50"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.839"
$ws.Range("E15").Value = "  -3.62%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.237.45"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "14.02"
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("D18").Value = "43.663.05"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "12.85"
$ws.Range("E19").Value = "  -9.55%  "
$ws.Range("D20").Value = "0.0₃0961"
$ws.Range("E20").Value = "  -4.22%  "
$ws.Range("D21").Value = "6.35"
$ws.Range("E21").Value = "  -4.45%  "
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("D23").Value = "3.05"
$ws.Range("E23").Value = "  -3.51%  "
$ws.Range("D24").Value = "233.41"
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("E25").Value = "  -7.58%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "10.26"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("D29").Value = "36.89"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").Value = "6.00"
$ws.Range("E30").Value = "  -6.71%  "
$ws.Range("D31").Value = "158.05"
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("D33").Value = "0.0832"
$ws.Range("E33").Value = "  -5.45%  "
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").Value = "3.20"
$ws.Range("E35").Value = "  -2.03%  "
$ws.Range("E36").Value = "  +4.93%  "
$ws.Range("D37").Value = "1.89"
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("E38").Value = "  -3.26%  "
$ws.Range("D39").Value = "16.23"
$ws.Range("E39").Value = "  +6.71%  "
$ws.Range("E40").Value = "  -4.27%  "
$ws.Range("D41").Value = "4.06"
$ws.Range("E41").Value = "  -8.37%  "
$ws.Range("D42").Value = "0.0309"
$ws.Range("E42").Value = "  -4.60%  "
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").Value = "1.732.93"
$ws.Range("E44").Value = "  -4.29%  "
$ws.Range("E45").Value = "  -5.63%  "
$ws.Range("B46").Value = "ordi"
$ws.Range("C46").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D46").Value = "73.67"
$ws.Range("E46").Value = "  -2.02%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "80.47"
$ws.Range("E47").Value = "  -4.69%  "
$ws.Range("D48").Value = "5.10"
$ws.Range("E48").Value = "  -4.38%  "
$ws.Range("D49").Value = "1.66"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").Value = "101.57"
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("D51").Value = "56.56"
$ws.Range("E51").Value = "  -4.30%  "
